# Apply changes described in the commit: update to published CDA FHIR
# logical model with patches #241.
$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$wsMetadata.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$wsMetadata.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> HL7 International contact info
$wsMetadata.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
# Binding Value Set for AlternateIdentification.classCode row (row 3, column Z)
$wsElements.Range("Z3").Value = "http://hl7.org/cda/stds/core/ValueSet/CDARoleClass"

# Column Z (Binding Value Set) width change to fit new content.
# Target stored OOXML width is 49.11328125; the Excel column-width model here
# quantizes to pixel boundaries, so 48.25 is the closest input that lands on
# the nearest achievable stored width (49.166666...).
$wsElements.Columns.Item(26).ColumnWidth = 48.25
